# Updates paises.xlsx (COVID-19 "Pais" sheet) with a refreshed data pull:
#  - the "last updated" timestamp in A1
#  - case numbers (Casos totales/Nuevos/activos/Recuperados/criticos/Muertes) for
#    several countries
#  - since the sheet is kept sorted descending by "Casos totales" (col B), a few
#    countries leapfrog their neighbours with the new numbers; those rows are
#    rewritten in place (country name + stats) rather than physically moved so the
#    row/record the rest of the sheet depends on (r="..") never shifts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Datos actualizados" timestamp -----------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 9 de Julio de 2020 a las 11:51"

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 3159514
$ws.Cells.Item(4, 3).Value = 582
$ws.Cells.Item(4, 4).Value = 1393256
$ws.Cells.Item(4, 5).Value = 1631385
$ws.Cells.Item(4, 7).Value = 11
$ws.Cells.Item(4, 8).Value = 134873

# Row 6: India
$ws.Cells.Item(6, 2).Value = 771129
$ws.Cells.Item(6, 3).Value = 2077
$ws.Cells.Item(6, 4).Value = 477685
$ws.Cells.Item(6, 5).Value = 272270
$ws.Cells.Item(6, 7).Value = 30
$ws.Cells.Item(6, 8).Value = 21174

# Row 20: Banglades
$ws.Cells.Item(20, 2).Value = 175494
$ws.Cells.Item(20, 3).Value = 3360
$ws.Cells.Item(20, 4).Value = 84544
$ws.Cells.Item(20, 5).Value = 88712
$ws.Cells.Item(20, 7).Value = 41
$ws.Cells.Item(20, 8).Value = 2238

# Row 29: Indonesia
$ws.Cells.Item(29, 2).Value = 70736
$ws.Cells.Item(29, 3).Value = 2657
$ws.Cells.Item(29, 4).Value = 32651
$ws.Cells.Item(29, 5).Value = 34668
$ws.Cells.Item(29, 7).Value = 58
$ws.Cells.Item(29, 8).Value = 3417

# Row 35: Kazajistan
$ws.Cells.Item(35, 4).Value = 35137
$ws.Cells.Item(35, 5).Value = 17620

# --- Ucrania / Paises Bajos / Filipinas / Oman re-rank (rows 37-40) ---------
# Row 37: Filipinas
$ws.Cells.Item(37, 1).Value = "Filipinas"
$ws.Cells.Item(37, 2).Value = 51754
$ws.Cells.Item(37, 3).Value = 1395
$ws.Cells.Item(37, 4).Value = 12813
$ws.Cells.Item(37, 5).Value = 37627
$ws.Cells.Item(37, 7).Value = 0
$ws.Cells.Item(37, 8).Value = 1314

# Row 38: Oman
$ws.Cells.Item(38, 1).Value = "Oman"
$ws.Cells.Item(38, 2).Value = 51725
$ws.Cells.Item(38, 3).Value = 1518
$ws.Cells.Item(38, 4).Value = 33021
$ws.Cells.Item(38, 5).Value = 18468
$ws.Cells.Item(38, 7).Value = 3
$ws.Cells.Item(38, 8).Value = 236

# Row 39: Ucrania
$ws.Cells.Item(39, 1).Value = "Ucrania"
$ws.Cells.Item(39, 2).Value = 51224
$ws.Cells.Item(39, 3).Value = 810
$ws.Cells.Item(39, 4).Value = 23784
$ws.Cells.Item(39, 5).Value = 26113
$ws.Cells.Item(39, 7).Value = 21
$ws.Cells.Item(39, 8).Value = 1327

# Row 40: Paises Bajos
$ws.Cells.Item(40, 1).Value = "Paises Bajos"
$ws.Cells.Item(40, 2).Value = 50746
$ws.Cells.Item(40, 4).Value = 0
$ws.Cells.Item(40, 5).Value = 0
$ws.Cells.Item(40, 8).Value = 6135

# Row 46: Polonia
$ws.Cells.Item(46, 2).Value = 36951
$ws.Cells.Item(46, 3).Value = 262
$ws.Cells.Item(46, 5).Value = 9923
$ws.Cells.Item(46, 7).Value = 9
$ws.Cells.Item(46, 8).Value = 1551

# --- Camerun / Marruecos re-rank (rows 65-66) --------------------------------
# Row 65: Marruecos
$ws.Cells.Item(65, 1).Value = "Marruecos"
$ws.Cells.Item(65, 2).Value = 14949
$ws.Cells.Item(65, 3).Value = 178
$ws.Cells.Item(65, 4).Value = 11372
$ws.Cells.Item(65, 5).Value = 3335
$ws.Cells.Item(65, 8).Value = 242

# Row 66: Camerun
$ws.Cells.Item(66, 1).Value = "Camerun"
$ws.Cells.Item(66, 2).Value = 14916
$ws.Cells.Item(66, 4).Value = 11525
$ws.Cells.Item(66, 5).Value = 3032
$ws.Cells.Item(66, 8).Value = 359

# Row 74: Noruega
$ws.Cells.Item(74, 2).Value = 8954
$ws.Cells.Item(74, 3).Value = 4
$ws.Cells.Item(74, 5).Value = 565

# Row 77: Malasia
$ws.Cells.Item(77, 2).Value = 8683
$ws.Cells.Item(77, 3).Value = 6
$ws.Cells.Item(77, 4).Value = 8499
$ws.Cells.Item(77, 5).Value = 63

# Row 83: Finlandia
$ws.Cells.Item(83, 2).Value = 7273
$ws.Cells.Item(83, 3).Value = 8
$ws.Cells.Item(83, 5).Value = 144

# Row 119: Lituania
$ws.Cells.Item(119, 2).Value = 1857
$ws.Cells.Item(119, 3).Value = 3
$ws.Cells.Item(119, 4).Value = 1564
$ws.Cells.Item(119, 5).Value = 214

# --- Barbados / Seychelles / Lesoto re-rank (rows 184-185) -------------------
# Row 184: Lesoto
$ws.Cells.Item(184, 1).Value = "Lesoto"

# Row 185: Seychelles
$ws.Cells.Item(185, 1).Value = "Seychelles"

# --- San Cristobal y Nieves / Islas Malvinas / Groenlandia re-rank (209-210) -
# Row 209: Groenlandia
$ws.Cells.Item(209, 1).Value = "Groenlandia"

# Row 210: Islas Malvinas
$ws.Cells.Item(210, 1).Value = "Islas Malvinas"
